$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell B1 keeps its "Percentage" text but its number format changes
# from 0.00% (numFmtId 10) to 0.00 (numFmtId 2).
$ws.Range("B1").NumberFormat = "0.00"

# Column B's default/background style picks up the new 0.00 number format
# (this also stamps every currently populated cell in column B).
$ws.Columns("B").NumberFormat = "0.00"

# B2 / B3: values re-cast from double to BigDecimal; restore their original
# percent display (reuses the still-existing 0.00% style).
$ws.Range("B2").Value = 0.0025000000000000001
$ws.Range("B2").NumberFormat = "0.00%"
$ws.Range("B3").Value = 0.00124
$ws.Range("B3").NumberFormat = "0.00%"

# B4's value is unchanged; also restore its percent display.
$ws.Range("B4").NumberFormat = "0.00%"

# New (empty) touched cell further out on the row, extending the used range.
$ws.Range("I7").NumberFormat = "General"

# B5 / B6 / B7: values changed and re-formatted as whole-number percentages.
$ws.Range("B5").Value = 0.01
$ws.Range("B5").NumberFormat = "0%"
$ws.Range("B6").Value = 0.01
$ws.Range("B6").NumberFormat = "0%"
$ws.Range("B7").Value = 2
$ws.Range("B7").NumberFormat = "0%"

# Final selection left on B6.
$ws.Range("B6").Select()
